$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 15:36"

# Row 4
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7047643
$ws.Range("C4").Value = 1427
$ws.Range("D4").Value = 4300813
$ws.Range("E4").Value = 2542253
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 71
$ws.Range("H4").Value = 204577

# Row 19
$ws.Range("A19").Value = "Arabia Saudita"
$ws.Range("B19").Value = 330798
$ws.Range("C19").Value = 552
$ws.Range("D19").Value = 312684
$ws.Range("E19").Value = 13572
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 4542

# Row 25
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 276061
$ws.Range("C25").Value = 510
$ws.Range("D25").Value = 246300
$ws.Range("E25").Value = 20276
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 9485

# Row 32
$ws.Range("A32").Value = "Catar"
$ws.Range("B32").Value = 123917
$ws.Range("C32").Value = 313
$ws.Range("D32").Value = 120766
$ws.Range("E32").Value = 2940
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 211

# Row 41
$ws.Range("A41").Value = "Paises Bajos"
$ws.Range("B41").Value = 98240
$ws.Range("C41").Value = 2245
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 10
$ws.Range("H41").Value = 6291

# Row 54
$ws.Range("A54").Value = "Nepal"
$ws.Range("B54").Value = 66632
$ws.Range("C54").Value = 1356
$ws.Range("D54").Value = 48061
$ws.Range("E54").Value = 18142
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 429

# Row 55
$ws.Range("A55").Value = "Barein"
$ws.Range("B55").Value = 65752
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 58626
$ws.Range("E55").Value = 6900
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 226

# Row 56
$ws.Range("A56").Value = "Costa Rica"
$ws.Range("B56").Value = 65602
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 25127
$ws.Range("E56").Value = 39730
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 745

# Row 61
$ws.Range("A61").Value = "Suiza"
$ws.Range("B61").Value = 50664
$ws.Range("C61").Value = 286
$ws.Range("D61").Value = 41800
$ws.Range("E61").Value = 6810
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 2054

# Row 74
$ws.Range("A74").Value = "Serbia"
$ws.Range("B74").Value = 32999
$ws.Range("C74").Value = 61
$ws.Range("D74").Value = 31536
$ws.Range("E74").Value = 720
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 743

# Row 79
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 25737
$ws.Range("C79").Value = 216
$ws.Range("D79").Value = 18359
$ws.Range("E79").Value = 6600
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 8
$ws.Range("H79").Value = 778

# Row 102
$ws.Range("A102").Value = "Tayikistan"
$ws.Range("B102").Value = 9432
$ws.Range("C102").Value = 44
$ws.Range("D102").Value = 8199
$ws.Range("E102").Value = 1160
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 73

# Row 112
$ws.Range("A112").Value = "Uganda"
$ws.Range("B112").Value = 6712
$ws.Range("C112").Value = 244
$ws.Range("D112").Value = 2778
$ws.Range("E112").Value = 3870
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 64

# Row 113
$ws.Range("A113").Value = "Birmania"
$ws.Range("B113").Value = 6471
$ws.Range("C113").Value = 320
$ws.Range("D113").Value = 1445
$ws.Range("E113").Value = 4926
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 100

# Row 131
$ws.Range("A131").Value = "Georgia"
$ws.Range("B131").Value = 3913
$ws.Range("C131").Value = 218
$ws.Range("D131").Value = 1574
$ws.Range("E131").Value = 2316
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 3
$ws.Range("H131").Value = 23

# Row 148
$ws.Range("A148").Value = "Islandia"
$ws.Range("B148").Value = 2419
$ws.Range("C148").Value = 42
$ws.Range("D148").Value = 2130
$ws.Range("E148").Value = 279
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 10

# Row 149
$ws.Range("A149").Value = "Guyana"
$ws.Range("B149").Value = 2402
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 1359
$ws.Range("E149").Value = 978
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 65

# Row 168
$ws.Range("A168").Value = "Vietnam"
$ws.Range("B168").Value = 1068
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 957
$ws.Range("E168").Value = 76
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 35

# Row 203
$ws.Range("A203").Value = "Fiyi"
$ws.Range("B203").Value = 32
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 28
$ws.Range("E203").Value = 2
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 2

# Row 214
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

